# Fruta / hortaliza, semanal
#
# A new weekly price-report record is inserted as row 232 of the data
# sheet (pushing the existing rows 232..331 down to 233..332, so the
# sheet grows from 331 to 332 data rows). Insert a full blank row first
# (shifting everything below it down), then populate the new row with
# the reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 232, shifting rows 232:331 down to 233:332.
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new record's data.
$ws.Cells.Item(232, 1).Value  = 9
$ws.Cells.Item(232, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(232, 3).Value  = "Metropolitana"
$ws.Cells.Item(232, 4).Value  = 45205
$ws.Cells.Item(232, 5).Value  = 13
$ws.Cells.Item(232, 6).Value  = "Fruta"
$ws.Cells.Item(232, 7).Value  = 100101
$ws.Cells.Item(232, 8).Value  = "Berries"
$ws.Cells.Item(232, 9).Value  = 100101001
$ws.Cells.Item(232, 10).Value = "Arándano (blue)"
$ws.Cells.Item(232, 11).Value = "Sin especificar"
$ws.Cells.Item(232, 12).Value = "Primera"
$ws.Cells.Item(232, 13).Value = 280
$ws.Cells.Item(232, 14).Value = 10000
$ws.Cells.Item(232, 15).Value = 11000
$ws.Cells.Item(232, 16).Value = 10536
$ws.Cells.Item(232, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(232, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(232, 19).Value = 7024
$ws.Cells.Item(232, 20).Value = 1.5

# Keep the date column's number format consistent with the rest of the
# column (inherited from the row above after the insert, but set it
# explicitly to be safe).
$ws.Cells.Item(232, 4).NumberFormat = $ws.Cells.Item(233, 4).NumberFormat
